$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Column width changes
# ---------------------------------------------------------------------------
$wsOverview.Range("E1").ColumnWidth = 29.9777047293527
$wsOverview.Range("F1").ColumnWidth = 29.9777047293527

$wsZhCn.Range("C1").ColumnWidth = 29.9777047293527
$wsZhCn.Range("I1").ColumnWidth = 40
$wsZhCn.Range("J1").ColumnWidth = 40

$wsDeDe.Range("C1").ColumnWidth = 29.9777047293527
$wsDeDe.Range("I1").ColumnWidth = 40
$wsDeDe.Range("J1").ColumnWidth = 40

# ---------------------------------------------------------------------------
# 3. zh-cn sheet: fill in "Latest Target File" (I) / "Latest Handback File" (J)
#    and the handback datetime (K), plus hyperlinks on I2/I3.
# ---------------------------------------------------------------------------
$zhCnFile1 = "04ea740f-b09f-4b70-af11-4b65e62206ba.md"
$zhCnFile2 = "5ad28dad-3920-4125-bac1-dd1defa03890.md"
$zhCnUrl1  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a5183b40492ef0ceafda257f1fb4c75a6ce7c300/e2e/04ea740f-b09f-4b70-af11-4b65e62206ba.md"
$zhCnUrl2  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a5183b40492ef0ceafda257f1fb4c75a6ce7c300/e2e/5ad28dad-3920-4125-bac1-dd1defa03890.md"

$wsZhCn.Range("I2").Value = $zhCnFile1
$wsZhCn.Range("J2").Value = "04ea740f-b09f-4b70-af11-4b65e62206ba.0f71a16a6badb4b329dda2489d1a45f25fee972c.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-31 15:01:25"

$wsZhCn.Range("I3").Value = $zhCnFile2
$wsZhCn.Range("J3").Value = "5ad28dad-3920-4125-bac1-dd1defa03890.4d8250def4d380e1cd26a8e21bf2f9635133db01.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-31 15:01:25"

# Recreate the hyperlinks list in row/column order (A2, I2, A3, I3) so the
# rIds for the new "Latest Target File" hyperlinks land right after their
# matching "File Name" hyperlinks, same as Excel does when it rewrites the
# sheet.
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhCnUrl1, "", "", $zhCnFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $zhCnUrl1, "", "", $zhCnFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $zhCnUrl2, "", "", $zhCnFile2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $zhCnUrl2, "", "", $zhCnFile2)

# ---------------------------------------------------------------------------
# 4. de-de sheet: same shape of change, different handback datetime.
# ---------------------------------------------------------------------------
$deDeFile1 = "04ea740f-b09f-4b70-af11-4b65e62206ba.md"
$deDeFile2 = "5ad28dad-3920-4125-bac1-dd1defa03890.md"
$deDeUrl1  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a5183b40492ef0ceafda257f1fb4c75a6ce7c300/e2e/04ea740f-b09f-4b70-af11-4b65e62206ba.md"
$deDeUrl2  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a5183b40492ef0ceafda257f1fb4c75a6ce7c300/e2e/5ad28dad-3920-4125-bac1-dd1defa03890.md"

$wsDeDe.Range("I2").Value = $deDeFile1
$wsDeDe.Range("J2").Value = "04ea740f-b09f-4b70-af11-4b65e62206ba.0f71a16a6badb4b329dda2489d1a45f25fee972c.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-31 15:01:40"

$wsDeDe.Range("I3").Value = $deDeFile2
$wsDeDe.Range("J3").Value = "5ad28dad-3920-4125-bac1-dd1defa03890.4d8250def4d380e1cd26a8e21bf2f9635133db01.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-31 15:01:40"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $deDeUrl1, "", "", $deDeFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $deDeUrl1, "", "", $deDeFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $deDeUrl2, "", "", $deDeFile2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $deDeUrl2, "", "", $deDeFile2)

Write-Host "Report for handback generated."
